$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw lat/lon readings (B2:B5) with new logged values
$ws.Range("B2").Value = 28.628305000000001
$ws.Range("B3").Value = -81.199619999999996
$ws.Range("B4").Value = 28.628391000000001
$ws.Range("B5").Value = -81.200012999999998

# Fix the bearing formula in D7: swap the ATAN2 arguments (y, x) order
$ws.Range("D7").Formula = "=ATAN2( SIN(C5-C3) * COS(C4), COS(C2) * SIN(C4) - SIN(C2) * COS(C4) * COS(C5-C3))"
